# Apply updated cryptocurrency price/volume values scraped on Sun May  7 10:03:38 UTC 2023.
# Cells hold text (not numeric) values in the source data, so a leading apostrophe is used
# to force text entry and the style is reset to Normal afterwards so no numeric formatting
# or quote-prefix styling is introduced on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.907.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.47%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.910.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.45%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = "'  -0.12%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.72%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.47%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.26%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.61%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -2.21%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.929.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.25%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'6.938"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.96%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.657"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.71%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.07058"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.28%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +0.00%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'83.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.21%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000009453"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.55%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'16.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.42%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'28.905.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.60%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.83%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -1.30%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.097"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.10%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'158.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.74%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'19.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.64%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'5.648"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.35%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'117.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.30%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.868"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.65%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.09303"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.39%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.8667"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.07%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'5.072"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.83%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -4.27%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.076"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.59%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.05716"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.67%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.165"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.03%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +0.10%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.02046"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.46%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'7.416"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.11%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5488"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.13%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1751"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.42%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.872"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'9.314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.84%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.5177"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.99%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -1.79%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.06904"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.58%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'2.092"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.04%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -1.92%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'110.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.80%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.000002555"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -11.71%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.2864"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.98%  "
$ws.Range("E51").Style = "Normal"
